# Append new Lancers listings snapshot (2026-01-25 12:38 JST).
#
# The sheet holds one "scrape" per append: row 2 is always the most recent
# run. This edit inserts a brand-new row 4 (AI medical chatbot listing),
# pushes the previously-seen rows 4-6 down to 5-7, keeps two of the old
# listings (rows that already existed) and appends two more brand-new
# listings at rows 8-9. Every row's timestamp (column A) is refreshed to
# the new run time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$timestamp = "2026-01-25 12:38:21"

# --- 1. Clear all existing hyperlink relationships up front ----------------
# (Individual Hyperlink.Delete() calls are unreliable in this host; the
# bulk Range.Hyperlinks.Delete() reliably clears every hyperlink on the
# sheet without touching cell values/styles, so we rebuild all of them
# below in final row order -- giving a clean rId1..rId8 sequence that
# lines up with the final F2:F9 URLs.)
$ws.Range("A1").Hyperlinks.Delete()

# --- 2. Write every data row's final values ---------------------------------

# Row 2 (unchanged content, refreshed timestamp)
$ws.Range("A2").Value = $timestamp
$ws.Range("B2").Value = "【Power Platform】請求書自動作成ツール開発(Automate / 外部API連携)"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5478394"
$ws.Range("G2").Value = 320
$ws.Range("H2").Value = "🔥API ◆ツール,開発"

# Row 3 (unchanged content, refreshed timestamp)
$ws.Range("A3").Value = $timestamp
$ws.Range("B3").Value = "【初心者・未経験OK】 AIを学びながら在宅で働くお仕事|月15〜20万円可能|スマホOK"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5478263"
$ws.Range("G3").Value = 303
$ws.Range("H3").Value = "🔥AI,Ai"

# Row 4 (brand-new listing)
$ws.Range("A4").Value = $timestamp
$ws.Range("B4").Value = "【音声収録】AI医療チャットボットと会話するだけ/津軽弁:方言話者募集"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "1,000 ~ 5,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5478603"
$ws.Range("G4").Value = 295
$ws.Range("H4").Value = "🔥AI,Ai"

# Row 5 (was row 4)
$ws.Range("A5").Value = $timestamp
$ws.Range("B5").Value = "「飲み会調整・店舗共有・終電管理・近距離マッチングを備えた飲み会支援アプリの開発依頼」"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5478300"
$ws.Range("G5").Value = 135
$ws.Range("H5").Value = "◆開発 ◇アプリ"

# Row 6 (was row 5)
$ws.Range("A6").Value = $timestamp
$ws.Range("B6").Value = "【Power Platform】外部業者連携システム開発(Apps / Automate)"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5478393"
$ws.Range("G6").Value = 125
$ws.Range("H6").Value = "◆開発,システム開発"

# Row 7 (brand-new listing)
$ws.Range("A7").Value = $timestamp
$ws.Range("B7").Value = "【共同開発メンバー募集】猫×IT|ブリーダー直販マッチングサービス「ねこ結び」"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5478555"
$ws.Range("G7").Value = 75
$ws.Range("H7").Value = "◆開発"

# Row 8 (was row 6)
$ws.Range("A8").Value = $timestamp
$ws.Range("B8").Value = "【Dify】LLMワークフロー開発パートナー募集 / 時間単価制"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E8").Value = "期限情報なし"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5478398"
$ws.Range("G8").Value = 68
$ws.Range("H8").Value = "◆開発"

# Row 9 (brand-new listing; no skill-summary text, same as source data)
$ws.Range("A9").Value = $timestamp
$ws.Range("B9").Value = "【緊急】スポーツ動画のボール追跡ロジック統合&iOS最適化(CoreML対応)"
$ws.Range("C9").Value = "システム開発"
$ws.Range("D9").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E9").Value = "期限情報なし"
$ws.Range("F9").Value = "https://www.lancers.jp/work/detail/5478514"
$ws.Range("G9").Value = 18

# --- 3. Rebuild hyperlinks on F2:F9 in row order ---------------------------
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5478394")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5478263")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5478603")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5478300")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5478393")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5478555")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5478398")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5478514")

# --- 4. Re-apply the Hyperlink cell style (Hyperlinks.Add() stamps its own
#        fresh style slot; putting the standard named "Hyperlink" style back
#        on each cell matches the pre-existing F2:F6 formatting). ----------
$ws.Range("F2").Style = "Hyperlink"
$ws.Range("F3").Style = "Hyperlink"
$ws.Range("F4").Style = "Hyperlink"
$ws.Range("F5").Style = "Hyperlink"
$ws.Range("F6").Style = "Hyperlink"
$ws.Range("F7").Style = "Hyperlink"
$ws.Range("F8").Style = "Hyperlink"
$ws.Range("F9").Style = "Hyperlink"
